# "added findaway normal for finance report"
#
# The per-file/per-currency rows (A2:E23) need to be re-sorted into
# ascending alphabetical order by the report key in column A
# (e.g. "87811004_0222_AU", "87811004_0222_BG", ... "87811004_0222_US"),
# keeping every row's B/C/D/E cells together with its A cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:E23")
$sortKey   = $ws.Range("A2:A23")

$dataRange.Sort($sortKey, 1)
